$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("105").Insert()

$ws.Range("A105").Value = 8
$ws.Range("B105").Value = "Terminal La Palmera de La Serena"
$ws.Range("C105").Value = "Coquimbo"
$ws.Range("D105").Value = 44601
$ws.Range("D105").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E105").Value = 4
$ws.Range("F105").Value = 100112003
$ws.Range("G105").Value = "Ajo"
$ws.Range("H105").Value = "Chino"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 560
$ws.Range("K105").Value = 18500
$ws.Range("L105").Value = 19000
$ws.Range("M105").Value = 18750
$ws.Range("N105").Value = "`$/caja 10 kilos"
$ws.Range("O105").Value = "China"
$ws.Range("P105").Value = 1875
$ws.Range("Q105").Value = 10
$ws.Range("R105").Value = "Hortaliza"
